$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 297
$ws1.Range("F3").Value = 1153
$ws1.Range("F4").Value = 2609

# Sheet "全部类型" (sheet4): update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 297
$ws4.Range("F5").Value = 1153
$ws4.Range("F6").Value = 2609
